$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the education-level codes in column A (rows 2-17): strip the
# "ZZ_" prefix per ISTAT's updated controlled vocabulary (e.g. "ZZ_NED" -> "NED").
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = [string]$cell.Text
    $cell.Value = $current -replace '^ZZ_', ''
}

# Move the view/selection to reflect where editing left off.
$ws.Range("A18").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
